$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly fruit/vegetable rows got re-ordered: what used to be rows 2/3
# (Quillota, week of 44915) now belongs in rows 4/5, and what used to be
# rows 4/5 (O'Higgins, week of 44911) now belongs in rows 2/3. Columns
# A,B,C,E,F,G,H,I,J,K,Q,T are identical between the swapped rows, so only
# D,L,M,N,O,P,R,S need to move.

$cols = @("D","L","M","N","O","P","R","S")

foreach ($col in $cols) {
    $topRef = $col + "2"
    $botRef = $col + "4"
    $topVal = $ws.Range($topRef).Value2
    $botVal = $ws.Range($botRef).Value2
    $ws.Range($topRef).Value2 = $botVal
    $ws.Range($botRef).Value2 = $topVal
}

foreach ($col in $cols) {
    $topRef = $col + "3"
    $botRef = $col + "5"
    $topVal = $ws.Range($topRef).Value2
    $botVal = $ws.Range($botRef).Value2
    $ws.Range($topRef).Value2 = $botVal
    $ws.Range($botRef).Value2 = $topVal
}
